$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.676693797111511
$ws.Range("B1").Value = 2.409196376800537
$ws.Range("C1").Value = 4.78117847442627
$ws.Range("D1").Value = 4.396365642547607
$ws.Range("E1").Value = 1.371733784675598
